$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: "We can cycle over all participants in a study and determine..."
#   -> "We can cycle over all participants in a study now and determine..."
# Insert "now " right after "in a study " (before "and determine").
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("We can cycle over all participants in a study ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Collapse(0)
$rng.InsertAfter("now ")

# ---------------------------------------------------------------------------
# Change 2: insert a new bullet paragraph right after "Received the sleep
# analysis and fell-asleep, wake-up algorithms from CamNTech" (and before the
# blank paragraph / "What is going to be done next:" heading), carrying the
# bookmark that used to sit at the end of the document's last bullet.
# ---------------------------------------------------------------------------
$rng2 = $d.Content
$rng2.Find.Execute("Received the sleep analysis and fell-asleep, wake-up algorithms from CamNTech", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$recvPara = $rng2.Paragraphs(1)
$recvPara.Range.InsertParagraphAfter() | Out-Null
$newPara = $recvPara.Next()
$newPara.Range.Text = "Have program to compare got-up and light’s-out times generated by the program and those found by people following the Motion Watch 8 Protocol (find in sleep analysis sheets for studies) to validate correctness and compare and contrast accuracy"

# ---------------------------------------------------------------------------
# Move the "_GoBack" bookmark from the end of the last bullet ("Demonstrate
# theoretical correctness...") to the end of the newly inserted bullet above.
# ---------------------------------------------------------------------------
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

$endRng = $d.Content
$endRng.Find.Execute("contrast accuracy", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$lastChar = $d.Range($endRng.End - 1, $endRng.End)
$d.Bookmarks.Add("_GoBack", $lastChar)

Write-Output "done"
